$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D hold text-formatted numeric-looking values (e.g. "305.01").
# Excel auto-converts a numeric-looking string assigned to .Value into a real
# number, which would lose the original text formatting (trailing zeros, etc).
# Forcing the cell to Text number-format ("@") before assignment keeps it literal,
# exactly like the source inlineStr cells.

$ws.Range("D2").Value = "41.913.00"
$ws.Range("E2").Value = "  -0.30%  "

$ws.Range("D3").Value = "2.263.30"
$ws.Range("E3").Value = "  -0.26%  "

$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "305.01"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "95.58"
$ws.Range("E6").Value = "  +2.59%  "

$ws.Range("E7").Value = "  -0.69%  "

$ws.Range("E9").Value = "  +0.50%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.06"
$ws.Range("E10").Value = "  +6.96%  "

$ws.Range("E11").Value = "  -1.06%  "

$ws.Range("E12").Value = "  -0.46%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.67"
$ws.Range("E13").Value = "  -0.29%  "

$ws.Range("D14").Value = "2.613.72"
$ws.Range("E14").Value = "  -0.29%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.37"
$ws.Range("E15").Value = "  +0.24%  "

$ws.Range("D16").Value = "2.257.16"
$ws.Range("E16").Value = "  -0.73%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.792"
$ws.Range("E17").Value = "  +0.85%  "

$ws.Range("D18").Value = "41.819.37"
$ws.Range("E18").Value = "  -0.22%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.37"
$ws.Range("E19").Value = "  -3.26%  "

$ws.Range("D20").Value = "0.0₃0901"

$ws.Range("E21").Value = "  -0.61%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "68.06"
$ws.Range("E22").Value = "  -0.18%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "237.11"
$ws.Range("E23").Value = "  -2.94%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.56"
$ws.Range("E24").Value = "  -1.41%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.93"
$ws.Range("E25").Value = "  -0.60%  "

$ws.Range("E26").Value = "  -0.01%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "23.65"
$ws.Range("E27").Value = "  -1.44%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "36.63"
$ws.Range("E28").Value = "  +4.60%  "

$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.12"
$ws.Range("E29").Value = "  +1.37%  "

$ws.Range("B30").Value = "Cosmos"
$ws.Range("C30").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.46"
$ws.Range("E30").Value = "  -2.29%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "160.35"
$ws.Range("E31").Value = "  +0.25%  "

$ws.Range("E32").Value = "  -2.35%  "

$ws.Range("E33").Value = "  +0.07%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.18"
$ws.Range("E34").Value = "  +4.85%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0735"
$ws.Range("E35").Value = "  -1.13%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "17.00"
$ws.Range("E36").Value = "  -0.85%  "

$ws.Range("E37").Value = "  -0.02%  "

$ws.Range("E38").Value = "  -0.67%  "

$ws.Range("E39").Value = "  +1.35%  "

$ws.Range("E40").Value = "  -2.33%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.00"
$ws.Range("E41").Value = "  +0.09%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.29"
$ws.Range("E42").Value = "  +2.24%  "

$ws.Range("D43").Value = "1.972.30"
$ws.Range("E43").Value = "  -2.20%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0282"
$ws.Range("E44").Value = "  -0.29%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "18.74"
$ws.Range("E45").Value = "  -5.21%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.92"
$ws.Range("E46").Value = "  +0.56%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.87"
$ws.Range("E47").Value = "  -3.85%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "52.98"
$ws.Range("E48").Value = "  -0.63%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "72.63"
$ws.Range("E49").Value = "  +0.21%  "

$ws.Range("E50").Value = "  -1.29%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "91.17"
$ws.Range("E51").Value = "  -0.67%  "
